$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.912.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "'1.648.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.03%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'308.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "'0.3885"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.90%  "

$ws.Range("D8").Value = "'0.3825"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "'51.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.08%  "

$ws.Range("D10").Value = "'1.349"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.15%  "

$ws.Range("D11").Value = "'1.000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.14%  "

$ws.Range("D12").Value = "'0.08434"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").Value = "'23.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("D14").Value = "'7.086"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").Value = "'7.775"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.12%  "

$ws.Range("D16").Value = "'0.00001309"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.67%  "

$ws.Range("D17").Value = "'1.646.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.13%  "

$ws.Range("D18").Value = "'94.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.54%  "

$ws.Range("D19").Value = "'0.06966"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("D20").Value = "'19.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.80%  "

$ws.Range("D21").Value = "'6.858"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.61%  "

$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "'13.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.22%  "

$ws.Range("D24").Value = "'23.909.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("D25").Value = "'2.483"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").Value = "'2.987"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.33%  "

$ws.Range("D27").Value = "'22.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.48%  "

$ws.Range("D28").Value = "'152.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.60%  "

$ws.Range("D29").Value = "'5.423"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.43%  "

$ws.Range("D30").Value = "'138.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").Value = "'7.728"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.78%  "

$ws.Range("D32").Value = "'2.487"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").Value = "'1.825.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.47%  "

$ws.Range("D34").Value = "'1.026"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.08%  "

$ws.Range("D35").Value = "'0.08024"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").Value = "'0.02948"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.49%  "

$ws.Range("D37").Value = "'6.690"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.93%  "

$ws.Range("D38").Value = "'10.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.35%  "

$ws.Range("D39").Value = "'0.2676"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.60%  "

$ws.Range("D40").Value = "'0.09094"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.06%  "

$ws.Range("D41").Value = "'0.7519"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "

$ws.Range("D42").Value = "'13.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("D43").Value = "'1.419"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("D44").Value = "'16.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.77%  "

$ws.Range("D45").Value = "'0.6902"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.71%  "

$ws.Range("D46").Value = "'2.436"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "

$ws.Range("D47").Value = "'4.069"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").Value = "'1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").Value = "'0.08281"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "

$ws.Range("D50").Value = "'134.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.89%  "

$ws.Range("D51").Value = "'1.221"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.35%  "
